# Fruta / hortaliza, semanal
# A new weekly observation (old row 2) is inserted at the top of the data
# table; every existing data row shifts down by one, and the last row's
# data is preserved as the new final row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 15
$newLastRow = 16

# 1) Read the current data rows (2..15) into memory, column by column,
#    before anything is overwritten.
$cols = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20)
$data = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $data[$r] = $rowVals
}

# 2) The newly created last row (16), column D, needs the same date-cell
#    style as the rest of column D. Apply the number format before
#    writing a value so the engine doesn't infer/allocate a throwaway
#    date style first.
$ws.Cells.Item($newLastRow, 4).NumberFormat = $ws.Cells.Item($lastRow, 4).NumberFormat()

# 3) Shift every existing row down by one: new row (r+1) gets the values
#    that used to live in row r. Walk bottom-up so we never clobber a
#    source row before it has been read (values already cached above, but
#    keep the safe order regardless).
for ($r = $lastRow; $r -ge 2; $r--) {
    $destRow = $r + 1
    $srcVals = $data[$r]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}

# 4) Write the brand-new observation into row 2.
$ws.Cells.Item(2, 1).Value = 8
$ws.Cells.Item(2, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(2, 3).Value = "Coquimbo"
$ws.Cells.Item(2, 4).Value = 44630
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100104
$ws.Cells.Item(2, 8).Value = "Frutos de pepita"
$ws.Cells.Item(2, 9).Value = 100104003
$ws.Cells.Item(2, 10).Value = "Membrillo"
$ws.Cells.Item(2, 11).Value = "Champion"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 16
$ws.Cells.Item(2, 14).Value = 400000
$ws.Cells.Item(2, 15).Value = 405000
$ws.Cells.Item(2, 16).Value = 402500
$ws.Cells.Item(2, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(2, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(2, 19).Value = 894
$ws.Cells.Item(2, 20).Value = 450
